# Applies the weekly update for "Fruta / hortaliza" data:
# - Duplicate the last data row (row 37) into a new row 38, preserving all its values.
# - Update row 37 with the new week's values (new date and new volume).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy existing row 37 (A:T) down into new row 38, preserving values/formatting.
$ws.Range("A37:T37").Copy($ws.Range("A38:T38"))

# Update row 37 with the new observation's values.
$ws.Range("D37").Value = 45239
$ws.Range("M37").Value = 100
